$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187, shifting existing rows 187-197 down to 188-198
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new "Ají" record
$ws.Cells.Item(187, 1).Value = 7
$ws.Cells.Item(187, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(187, 3).Value = "Ñuble"
$ws.Cells.Item(187, 4).Value = "2023-04-25"
$ws.Cells.Item(187, 5).Value = 16
$ws.Cells.Item(187, 6).Value = 100112021
$ws.Cells.Item(187, 7).Value = "Ají"
$ws.Cells.Item(187, 8).Value = "Cristal"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 50
$ws.Cells.Item(187, 11).Value = 17000
$ws.Cells.Item(187, 12).Value = 17000
$ws.Cells.Item(187, 13).Value = 17000
$ws.Cells.Item(187, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(187, 15).Value = "Región del Maule"
$ws.Cells.Item(187, 16).Value = 680
$ws.Cells.Item(187, 17).Value = 25
$ws.Cells.Item(187, 18).Value = "Hortaliza"
